# Deploy the implementation guide.
#
# The "Metadata" sheet's Property/Value table gains a new "Jurisdiction"
# row (inserted right after "Contact", pushing everything below it down
# by one row), and two of the existing property values are refreshed:
#   - Date    -> 2024-10-02T15:04:17+00:00
#   - Contact -> Ferlab.bio (http://example.org/example-publisher)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new blank row right below the "Contact" row (row 10), shifting
# rows 11-21 down to 12-22.
$ws.Rows("11:11").Insert()

# The freshly inserted row doesn't inherit the table's usual formatting;
# copy it over from the row directly above ("Contact") so the new row
# matches the rest of the property table.
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

# Populate the new "Jurisdiction" row (no value given).
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# Refresh the publication Date and Contact values.
$ws.Range("B8").Value = "2024-10-02T15:04:17+00:00"
$ws.Range("B10").Value = "Ferlab.bio (http://example.org/example-publisher)"
